$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-51 with refreshed crypto data.
# D column values are written as text (NumberFormat "@" + Style reset to "Normal" afterwards)
# so that numeric-looking strings (e.g. "1.000", "0.9998") are preserved exactly as text
# instead of being coerced into numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.345.86'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.92%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.625.56'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.10%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9998'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '304.21'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.29%  '

$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '51.86'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3618'
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08095'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.227'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.83%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.25%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.69'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.556'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001248'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.227'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.66%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.622.54'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.91%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.54'
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06917'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.54%  '

$ws.Range("E20").Value = '  -3.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.418'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.94%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '23.340.17'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.99%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.70'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.55%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.231'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.87%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.445'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.70%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.81%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '149.89'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.288'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.65%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.44'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.59%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.300'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.804.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.43%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.793'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.01'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9527'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.72%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02779'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.25%  '

$ws.Range("E37").Value = '  -1.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08813'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.64%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.082'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.83%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.07128'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.361'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.87%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7049'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.16'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.29'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.45%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6444'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.37%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.317'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.39%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9989'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.992'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.24%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07978'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.197'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.71%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '125.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.80%  '
